$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price and 1h volume-change columns); a handful of
# rows also have coin identity (name/link) changes because two rows swapped rank.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.813.05"
$ws.Range("E2").Value = "  -0.01%  "
# Row 3 - Ethereum
$ws.Range("D3").Value = "1.632.08"
$ws.Range("E3").Value = "  +0.20%  "
# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.26%  "
# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.59"
$ws.Range("E5").Value = "  -0.44%  "
# Row 6 - XRP
$ws.Range("E6").Value = "  -0.26%  "
# Row 7 - USDC
$ws.Range("E7").Value = "  -0.23%  "
# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.22%  "
# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0641"
$ws.Range("E9").Value = "  -0.65%  "
# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.86"
$ws.Range("E10").Value = "  +2.62%  "
# Row 11 - TRON
$ws.Range("E11").Value = "  +0.08%  "
# Row 12 - Polkadot
$ws.Range("E12").Value = "  -0.22%  "
# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.636.47"
$ws.Range("E13").Value = "  +0.51%  "
# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.859.19"
$ws.Range("E14").Value = "  +0.29%  "
# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.558"
$ws.Range("E15").Value = "  +0.03%  "
# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").Value = "  +1.80%  "
# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.95"
$ws.Range("E17").Value = "  -0.45%  "
# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.822.39"
$ws.Range("E18").Value = "  +0.00%  "
# Row 19 - Dai
$ws.Range("E19").Value = "  -0.28%  "
# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.60"
$ws.Range("E20").Value = "  +0.17%  "
# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  +1.98%  "
# Row 22 - Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.90"
$ws.Range("E22").Value = "  +0.99%  "
# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.21"
$ws.Range("E23").Value = "  +3.34%  "
# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  -0.25%  "
# Row 25 - Toncoin
$ws.Range("E25").Value = "  -3.36%  "
# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.48"
$ws.Range("E26").Value = "  -0.52%  "
# Row 27 - Stellar
$ws.Range("E27").Value = "  -3.50%  "
# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.81"
$ws.Range("E28").Value = "  +1.45%  "
# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("E29").Value = "  +0.77%  "
# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.10%  "
# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0494"
$ws.Range("E31").Value = "  +1.95%  "
# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  +0.40%  "
# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.88%  "
# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.98%  "
# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +0.16%  "
# Row 36 - ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.901"
$ws.Range("E36").Value = "  +0.84%  "
# Row 37 - MXToken
$ws.Range("E37").Value = "  +0.23%  "
# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +0.71%  "
# Row 39 - Maker
$ws.Range("D39").Value = "1.119.24"
$ws.Range("E39").Value = "  -1.19%  "
# Row 40 - VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0155"
$ws.Range("E40").Value = "  -0.08%  "
# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.03%  "
# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.52"
$ws.Range("E42").Value = "  -0.78%  "
# Row 43 - Quant
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.60"
$ws.Range("E43").Value = "  +2.37%  "
# Row 44 - TrustWalletToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.800"
$ws.Range("E44").Value = "  +0.58%  "
# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.0₆0108"
$ws.Range("E45").Value = "  -2.78%  "
# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.44"
$ws.Range("E46").Value = "  +0.86%  "
# Row 47 - Mantle
$ws.Range("E47").Value = "  -5.01%  "
# Row 48 - SynthetixNetwork
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +10.31%  "
# Row 49 - Cronos (was EnergySwap)
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  -0.56%  "
# Row 50 - EnergySwap (was Cronos)
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  +0.00%  "
# Row 51 - Frax
$ws.Range("E51").Value = "  +0.04%  "
